$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(860.6843515803669, 3.912201598092577, 90.86668033906227, 1.238472899808378, 863.1398315429688, 3.862115621566772, 90.77709197998047, 1.235077261924744)
    3 = @(513.1940128816041, 2.332700058552746, 79.44745434032257, 0.9745546493291284, 521.2578735351562, 2.324284553527832, 79.54720306396484, 0.9769288897514343)
    4 = @(804.1806264807286, 3.655366484003312, 65.6520626667446, 1.837266748320942, 802.1593017578125, 3.630875825881958, 65.55773162841797, 1.849378705024719)
    5 = @(664.7789928912285, 3.02172269496013, 83.12408038379573, 1.02995787400382, 666.68408203125, 2.972813367843628, 83.07846832275391, 1.023933291435242)
    6 = @(340.8504879206644, 1.549320399639383, 53.43999995745475, 0.8680711450312685, 340.1817932128906, 1.550617694854736, 53.47706604003906, 0.8605641722679138)
    7 = @(687.8861331246297, 3.126755150566499, 72.95857147479302, 1.414364753068805, 694.6416625976562, 3.135890483856201, 73.25852203369141, 1.41679584980011)
    8 = @(802.509261041841, 3.647769368372005, 55.43555057365677, 1.906398113915365, 795.4888305664062, 3.613675117492676, 55.22809600830078, 1.907754063606262)
    9 = @(531.5882751087696, 2.416310341403498, 83.54535462530259, 1.042156317201361, 545.5994873046875, 2.397069692611694, 82.85211181640625, 1.055340886116028)
    10 = @(665.0718122683322, 3.023053692128783, 62.66461948322814, 2.219185106925259, 661.8773193359375, 2.931782484054565, 62.35832977294922, 2.216866016387939)
    11 = @(441.7737518169947, 2.008062508259067, 68.67861429948152, 0.788055812582565, 443.8110961914062, 1.987303018569946, 68.86837768554688, 0.7809013724327087)
    12 = @(606.0441660017802, 2.754746209099001, 53.56798490560772, 2.237985935930394, 601.9494018554688, 2.734387159347534, 53.1558952331543, 2.23185133934021)
    13 = @(1081.269059136127, 4.914859359709668, 81.82801710489395, 1.746387887999544, 1086.56396484375, 4.909054756164551, 82.27091979980469, 1.757545471191406)
    14 = @(1061.932245949086, 4.826964754314026, 65.36674348599399, 2.42927699554501, 1049.491455078125, 4.788876533508301, 65.35366821289062, 2.441235303878784)
    15 = @(629.3744366956354, 2.86079289407107, 65.84493794374413, 1.816211615007142, 632.682861328125, 2.806559085845947, 65.55947113037109, 1.820563316345215)
    16 = @(1005.239664868703, 4.569271203948652, 57.67492178198188, 1.726740020102551, 1006.644104003906, 4.59144115447998, 58.16543579101562, 1.722370147705078)
    17 = @(650.3490175894341, 2.956131898133791, 85.32906898357351, 0.7628025077028043, 668.15576171875, 3.014429807662964, 85.2786865234375, 0.7689965963363647)
    18 = @(1058.553330603898, 4.811606048199535, 84.07738978894008, 1.26383921309223, 1055.292236328125, 4.83100414276123, 84.56607055664062, 1.256818175315857)
    19 = @(364.7597453313526, 1.657998842415239, 69.4877607209408, 0.89949135755155, 375.8826599121094, 1.641098022460938, 69.33924102783203, 0.8990448117256165)
    20 = @(832.7068392230001, 3.785031087377273, 84.69199126132031, 1.525891463878627, 839.0126342773438, 3.756873846054077, 84.62397003173828, 1.522988557815552)
    21 = @(381.9063300471682, 1.735937863850765, 71.77123256861205, 0.9262474468428752, 392.1371765136719, 1.715217351913452, 71.59212493896484, 0.9299800395965576)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
